# Replace the 25 multiplication problems in the table with new values,
# matching the commit "Update master to output generated at 4250d90".
$d = $word.ActiveDocument

$d.Content.Find.Execute('933×9=', $true, $false, $false, $false, $false, $true, 1, $false, '820×8=', 2) | Out-Null
$d.Content.Find.Execute('814×9=', $true, $false, $false, $false, $false, $true, 1, $false, '227×7=', 2) | Out-Null
$d.Content.Find.Execute('473×2=', $true, $false, $false, $false, $false, $true, 1, $false, '193×6=', 2) | Out-Null
$d.Content.Find.Execute('999×2=', $true, $false, $false, $false, $false, $true, 1, $false, '474×9=', 2) | Out-Null
$d.Content.Find.Execute('402×8=', $true, $false, $false, $false, $false, $true, 1, $false, '612×4=', 2) | Out-Null
$d.Content.Find.Execute('743×5=', $true, $false, $false, $false, $false, $true, 1, $false, '655×2=', 2) | Out-Null
$d.Content.Find.Execute('281×3=', $true, $false, $false, $false, $false, $true, 1, $false, '416×5=', 2) | Out-Null
$d.Content.Find.Execute('354×5=', $true, $false, $false, $false, $false, $true, 1, $false, '731×7=', 2) | Out-Null
$d.Content.Find.Execute('612×5=', $true, $false, $false, $false, $false, $true, 1, $false, '631×2=', 2) | Out-Null
$d.Content.Find.Execute('327×2=', $true, $false, $false, $false, $false, $true, 1, $false, '537×7=', 2) | Out-Null
$d.Content.Find.Execute('823×3=', $true, $false, $false, $false, $false, $true, 1, $false, '631×8=', 2) | Out-Null
$d.Content.Find.Execute('452×6=', $true, $false, $false, $false, $false, $true, 1, $false, '362×9=', 2) | Out-Null
$d.Content.Find.Execute('482×6=', $true, $false, $false, $false, $false, $true, 1, $false, '497×4=', 2) | Out-Null
$d.Content.Find.Execute('378×7=', $true, $false, $false, $false, $false, $true, 1, $false, '132×7=', 2) | Out-Null
$d.Content.Find.Execute('900×3=', $true, $false, $false, $false, $false, $true, 1, $false, '498×3=', 2) | Out-Null
$d.Content.Find.Execute('351×4=', $true, $false, $false, $false, $false, $true, 1, $false, '121×8=', 2) | Out-Null
$d.Content.Find.Execute('544×6=', $true, $false, $false, $false, $false, $true, 1, $false, '543×9=', 2) | Out-Null
$d.Content.Find.Execute('416×7=', $true, $false, $false, $false, $false, $true, 1, $false, '132×3=', 2) | Out-Null
$d.Content.Find.Execute('999×8=', $true, $false, $false, $false, $false, $true, 1, $false, '555×8=', 2) | Out-Null
$d.Content.Find.Execute('104×7=', $true, $false, $false, $false, $false, $true, 1, $false, '824×2=', 2) | Out-Null
$d.Content.Find.Execute('262×4=', $true, $false, $false, $false, $false, $true, 1, $false, '572×9=', 2) | Out-Null
$d.Content.Find.Execute('584×4=', $true, $false, $false, $false, $false, $true, 1, $false, '809×2=', 2) | Out-Null
$d.Content.Find.Execute('144×7=', $true, $false, $false, $false, $false, $true, 1, $false, '612×3=', 2) | Out-Null
$d.Content.Find.Execute('898×3=', $true, $false, $false, $false, $false, $true, 1, $false, '194×5=', 2) | Out-Null
$d.Content.Find.Execute('298×5=', $true, $false, $false, $false, $false, $true, 1, $false, '936×6=', 2) | Out-Null
